$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with swapped coin identity (B/C/D/E) ---
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "'3.482.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "'0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'5.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "'0.535"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'2.832.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'27.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "

# --- Rows with only price / volume updates ---
$ws.Range("D2").Value = "'67.186.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "'3.477.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'592.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'177.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +4.65%  "
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "'4.081.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "'31.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.93%  "
$ws.Range("D15").Value = "'0.135"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'67.249.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'3.479.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'6.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "'14.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "'387.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'7.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'74.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "'0.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'10.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'6.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'23.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "'7.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").Value = "'164.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").Value = "'1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "'2.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.81%  "
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").Value = "'4.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D46").Value = "'26.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "'0.0721"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").Value = "'41.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "'335.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  -2.37%  "
